$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "CT113"
$ws.Range("C13").Value = "ASFLI"
$ws.Range("D13").Value = "REGISTERPROJECT"
$ws.Range("E13").Value = "PENDING"
$ws.Range("F13").Value = 8
$ws.Range("G13").Style = "Normal"
$ws.Range("H13").Style = "Normal"
